$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 114 (existing fixture id=112 is now a played match -> new id/odds, and
# gains FTHG/FTAG/FTR (H/I/J) plus PL_AhOver/PL_AhUnder (AB/AC) results).
# ---------------------------------------------------------------------------
$ws.Range("B114").Value = 7749770
$ws.Range("E114").Value = 45381.35416666666
$ws.Range("F114").Value = "Bengaluru"
$ws.Range("G114").Value = "Odisha FC"
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = "D"
$ws.Range("K114").Value = 2.55
$ws.Range("L114").Value = 3.3
$ws.Range("M114").Value = 2.55
$ws.Range("N114").Value = 2.625
$ws.Range("O114").Value = 3.25
$ws.Range("P114").Value = 2.5
$ws.Range("Q114").Value = 0
$ws.Range("R114").Value = 2
$ws.Range("S114").Value = 1.85
$ws.Range("T114").Value = 2.5
$ws.Range("U114").Value = 1.975
$ws.Range("V114").Value = 1.825
$ws.Range("W114").Value = -1
$ws.Range("X114").Value = 2.25
$ws.Range("Y114").Value = -1
$ws.Range("AA114").Value = 0
$ws.Range("AB114").Value = -1
$ws.Range("AC114").Value = 0.825

# ---------------------------------------------------------------------------
# Row 115 (existing fixture id=113 is now a played match -> new id/odds, and
# gains FTHG/FTAG/FTR (H/I/J) plus PL columns W..AC).
# ---------------------------------------------------------------------------
$ws.Range("B115").Value = 7749469
$ws.Range("E115").Value = 45381.45833333334
$ws.Range("F115").Value = "Jamshedpur FC"
$ws.Range("G115").Value = "Kerala Blasters"
$ws.Range("H115").Value = 1
$ws.Range("I115").Value = 1
$ws.Range("J115").Value = "D"
$ws.Range("K115").Value = 2.2
$ws.Range("L115").Value = 3.25
$ws.Range("M115").Value = 3
$ws.Range("N115").Value = 2.15
$ws.Range("O115").Value = 3.3
$ws.Range("P115").Value = 3.1
$ws.Range("Q115").Value = -0.25
$ws.Range("T115").Value = 2.5
$ws.Range("U115").Value = 1.95
$ws.Range("V115").Value = 1.85
$ws.Range("W115").Value = -1
$ws.Range("X115").Value = 2.3
$ws.Range("Y115").Value = -1
$ws.Range("Z115").Value = -0.5
$ws.Range("AA115").Value = 0.5
$ws.Range("AB115").Value = -1
$ws.Range("AC115").Value = 0.8500000000000001

# ---------------------------------------------------------------------------
# Row 116 (brand-new future fixture, the old id=7749875 Hyderabad FC vs
# Mumbai City FC record moved down with refreshed odds; no result yet).
# ---------------------------------------------------------------------------
$a = $ws.Range("A116")
$a.Value = 114
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$a.Borders.LineStyle = 1
$a.Borders.Weight = 2

$ws.Range("B116").Value = 7749875
$ws.Range("C116").Value = "India Super League"
$ws.Range("D116").Value = "India Super League"

$e = $ws.Range("E116")
$e.Value = 45383.45833333334
$e.NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("F116").Value = "Hyderabad FC"
$ws.Range("G116").Value = "Mumbai City FC"
$ws.Range("K116").Value = 9.5
$ws.Range("L116").Value = 5.5
$ws.Range("M116").Value = 1.25
$ws.Range("N116").Value = 10
$ws.Range("O116").Value = 5.75
$ws.Range("P116").Value = 1.25
$ws.Range("Q116").Value = 1.75
$ws.Range("R116").Value = 1.9
$ws.Range("S116").Value = 1.9
$ws.Range("T116").Value = 3
$ws.Range("U116").Value = 1.9
$ws.Range("V116").Value = 1.9
$ws.Range("W116").Value = 0
$ws.Range("X116").Value = 0
$ws.Range("Y116").Value = 0
$ws.Range("Z116").Value = 0
$ws.Range("AA116").Value = 0
